$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (row 2 through 111) holds the "Förändrad" date serial 45189 (2023-09-20).
# Bump it by one day to 45190 (2023-09-21) for every data row.
$ws.Range("C2:C111").Value = 45190
